# Scheduled market-data refresh for Pandaemonium_Profits leve tables.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# per job/class sheet to the latest Universalis price snapshot.
$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 3483.6667
$ws.Range("I10").Value = 451
$ws.Range("J10").Value = 5000
$ws.Range("K10").Value = 451
$ws.Range("L10").Value = 5000
$ws.Range("M10").Value = -158
$ws.Range("N10").Value = -5586
$ws.Range("H16").Value = 10000
$ws.Range("I16").Value = 10000
$ws.Range("K16").Value = 10000
$ws.Range("M16").Value = -9770
$ws.Range("H33").Value = 330.25
$ws.Range("I33").Value = 321.92307
$ws.Range("J33").Value = 345.7143
$ws.Range("K33").Value = 321.92307
$ws.Range("L33").Value = 345.7143
$ws.Range("M33").Value = -92.92307
$ws.Range("N33").Value = -803.7143
$ws.Range("H43").Value = 2501
$ws.Range("I43").Value = 2000
$ws.Range("J43").Value = 3002
$ws.Range("K43").Value = 2000
$ws.Range("L43").Value = 3002
$ws.Range("M43").Value = -1931
$ws.Range("N43").Value = -3140
$ws.Range("H62").Value = 3637.125
$ws.Range("I62").Value = 2824.25
$ws.Range("J62").Value = 4450
$ws.Range("K62").Value = 2824.25
$ws.Range("L62").Value = 4450
$ws.Range("M62").Value = -2200.25
$ws.Range("N62").Value = -5698
$ws.Range("H65").Value = 3637.125
$ws.Range("I65").Value = 2824.25
$ws.Range("J65").Value = 4450
$ws.Range("K65").Value = 14121.25
$ws.Range("L65").Value = 22250
$ws.Range("M65").Value = -11001.25
$ws.Range("N65").Value = -28490
$ws.Range("H125").Value = 22574.4
$ws.Range("J125").Value = 25768
$ws.Range("L125").Value = 231912
$ws.Range("N125").Value = -236832
$ws.Range("H138").Value = 3793.7646
$ws.Range("J138").Value = 4203.544
$ws.Range("L138").Value = 12610.632
$ws.Range("N138").Value = -22890.632

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 35499
$ws.Range("J62").Value = 35499
$ws.Range("L62").Value = 35499
$ws.Range("N62").Value = -36747
$ws.Range("H65").Value = 35499
$ws.Range("J65").Value = 35499
$ws.Range("L65").Value = 106497
$ws.Range("N65").Value = -112737
$ws.Range("H74").Value = 5568.0557
$ws.Range("I74").Value = 3141.2307
$ws.Range("J74").Value = 11877.8
$ws.Range("K74").Value = 3141.2307
$ws.Range("L74").Value = 11877.8
$ws.Range("M74").Value = -2267.2307
$ws.Range("N74").Value = -13625.8
$ws.Range("H77").Value = 5568.0557
$ws.Range("I77").Value = 3141.2307
$ws.Range("J77").Value = 11877.8
$ws.Range("K77").Value = 15706.1535
$ws.Range("L77").Value = 59389
$ws.Range("M77").Value = -11338.1535
$ws.Range("N77").Value = -68125
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("M109").ClearContents()
$ws.Range("N109").ClearContents()

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 825
$ws.Range("J64").Value = 900
$ws.Range("L64").Value = 900
$ws.Range("N64").Value = -1350
$ws.Range("H67").Value = 825
$ws.Range("J67").Value = 900
$ws.Range("L67").Value = 900
$ws.Range("N67").Value = -2460
$ws.Range("H75").Value = 8736.444
$ws.Range("I75").Value = 5828.5
$ws.Range("J75").Value = 32000
$ws.Range("K75").Value = 5828.5
$ws.Range("L75").Value = 32000
$ws.Range("M75").Value = -4892.5
$ws.Range("N75").Value = -33872
$ws.Range("H78").Value = 8736.444
$ws.Range("I78").Value = 5828.5
$ws.Range("J78").Value = 32000
$ws.Range("K78").Value = 17485.5
$ws.Range("L78").Value = 96000
$ws.Range("M78").Value = -12805.5
$ws.Range("N78").Value = -105360
$ws.Range("H99").Value = 2627.3333
$ws.Range("I99").Value = 2571.4285
$ws.Range("J99").Value = 2676.25
$ws.Range("K99").Value = 2571.4285
$ws.Range("L99").Value = 2676.25
$ws.Range("M99").Value = -1073.4285
$ws.Range("N99").Value = -5672.25

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2337.7454
$ws.Range("I31").Value = 1652.4889
$ws.Range("J31").Value = 5421.4
$ws.Range("K31").Value = 1652.4889
$ws.Range("L31").Value = 5421.4
$ws.Range("M31").Value = -1357.4889
$ws.Range("N31").Value = -6011.4
$ws.Range("H34").Value = 2337.7454
$ws.Range("I34").Value = 1652.4889
$ws.Range("J34").Value = 5421.4
$ws.Range("K34").Value = 1652.4889
$ws.Range("L34").Value = 5421.4
$ws.Range("M34").Value = -1450.4889
$ws.Range("N34").Value = -5825.4
$ws.Range("H62").Value = 3188.875
$ws.Range("J62").Value = 3753
$ws.Range("L62").Value = 3753
$ws.Range("N62").Value = -5001
$ws.Range("H65").Value = 3188.875
$ws.Range("J65").Value = 3753
$ws.Range("L65").Value = 18765
$ws.Range("N65").Value = -25005
$ws.Range("H94").Value = 1606.375
$ws.Range("I94").Value = 1684.6
$ws.Range("J94").Value = 1476
$ws.Range("K94").Value = 1684.6
$ws.Range("L94").Value = 1476
$ws.Range("M94").Value = -1233.6
$ws.Range("N94").Value = -2378
$ws.Range("H134").Value = 4083.795
$ws.Range("I134").Value = 3191.182
$ws.Range("J134").Value = 4434.4644
$ws.Range("K134").Value = 9573.545999999998
$ws.Range("L134").Value = 13303.3932
$ws.Range("M134").Value = -7038.545999999998
$ws.Range("N134").Value = -18373.3932

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 19800.64
$ws.Range("I131").Value = 685.75
$ws.Range("J131").Value = 28795.883
$ws.Range("K131").Value = 2057.25
$ws.Range("L131").Value = 86387.649
$ws.Range("M131").Value = 2982.75
$ws.Range("N131").Value = -96467.649
$ws.Range("H136").Value = 2609.3635
$ws.Range("I136").Value = 1011.1111
$ws.Range("J136").Value = 3715.8462
$ws.Range("K136").Value = 3033.3333
$ws.Range("L136").Value = 11147.5386
$ws.Range("M136").Value = 2066.6667
$ws.Range("N136").Value = -21347.5386
$ws.Range("H137").Value = 22753042
$ws.Range("I137").Value = 41668070
$ws.Range("K137").Value = 125004210
$ws.Range("M137").Value = -124999110
$ws.Range("H139").Value = 1958634.9
$ws.Range("I139").Value = 3355444.8
$ws.Range("K139").Value = 10066334.4
$ws.Range("M139").Value = -10061194.4
$ws.Range("H141").Value = 2923.889
$ws.Range("I141").Value = 2512.6667
$ws.Range("K141").Value = 7538.000100000001
$ws.Range("M141").Value = -2358.000100000001

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 8600
$ws.Range("I18").Value = 7666.6665
$ws.Range("K18").Value = 7666.6665
$ws.Range("M18").Value = -7493.6665
$ws.Range("H20").Value = 10000
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H22").Value = 6904
$ws.Range("I22").Value = 712
$ws.Range("K22").Value = 712
$ws.Range("M22").Value = -419
